# Apply updated cryptocurrency price/volume data per upstream source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.534.05'
$ws.Range('E2').Value = '  -0.25%  '

# Row 3
$ws.Range('D3').Value = '3.669.33'
$ws.Range('E3').Value = '  -1.19%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '619.46'
$ws.Range('E5').Value = '  -7.86%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.11'
$ws.Range('E6').Value = '  -1.67%  '

# Row 7
$ws.Range('E7').Value = '  +0.11%  '

# Row 8
$ws.Range('E8').Value = '  -0.57%  '

# Row 9
$ws.Range('E9').Value = '  -2.27%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.16'
$ws.Range('E10').Value = '  +1.08%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.439'
$ws.Range('E11').Value = '  -1.35%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000228'
$ws.Range('E12').Value = '  -3.30%  '

# Row 13
$ws.Range('D13').Value = '4.288.63'
$ws.Range('E13').Value = '  -1.17%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.25'
$ws.Range('E14').Value = '  -1.99%  '

# Row 15
$ws.Range('D15').Value = '3.670.45'
$ws.Range('E15').Value = '  -0.36%  '

# Row 16
$ws.Range('D16').Value = '69.617.92'
$ws.Range('E16').Value = '  -0.18%  '

# Row 17
$ws.Range('E17').Value = '  +0.55%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.49'
$ws.Range('E18').Value = '  -0.50%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '15.83'
$ws.Range('E19').Value = '  -2.93%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.28'
$ws.Range('E20').Value = '  +4.50%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '468.91'
$ws.Range('E21').Value = '  -1.23%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.647'
$ws.Range('E22').Value = '  -1.33%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '79.55'
$ws.Range('E23').Value = '  -1.16%  '

# Row 24
$ws.Range('D24').Value = '3.817.51'
$ws.Range('E24').Value = '  -1.10%  '

# Row 25
$ws.Range('E25').Value = '  +0.03%  '

# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.04'
$ws.Range('E26').Value = '  +0.36%  '

# Row 27
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000122'
$ws.Range('E27').Value = '  -4.92%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.68'
$ws.Range('E28').Value = '  -4.93%  '

# Row 29
$ws.Range('E29').Value = '  -3.46%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.66'
$ws.Range('E30').Value = '  -4.42%  '

# Row 31
$ws.Range('E31').Value = '  +0.07%  '

# Row 32
$ws.Range('E32').Value = '  -2.11%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.53'
$ws.Range('E33').Value = '  -1.60%  '

# Row 34
$ws.Range('E34').Value = '  -2.52%  '

# Row 35
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = '3.670.65'
$ws.Range('E35').Value = '  -0.85%  '

# Row 36
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.37'
$ws.Range('E36').Value = '  -3.87%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '8.26'
$ws.Range('E37').Value = '  -3.57%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '178.47'
$ws.Range('E39').Value = '  +3.02%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.05%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.22'
$ws.Range('E41').Value = '  -1.43%  '

# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.76'
$ws.Range('E42').Value = '  -5.63%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0889'
$ws.Range('E43').Value = '  -2.75%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.924'
$ws.Range('E44').Value = '  -1.97%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '46.66'
$ws.Range('E45').Value = '  -0.97%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '28.95'
$ws.Range('E46').Value = '  +4.34%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.69'
$ws.Range('E47').Value = '  -3.27%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.83'
$ws.Range('E48').Value = '  -0.89%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.000261'
$ws.Range('E49').Value = '  -7.69%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.03'
$ws.Range('E50').Value = '  -5.02%  '

# Row 51
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.261'
$ws.Range('E51').Value = '  -2.77%  '
